$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 163: drop the trailing empty R163 cell (imagem) ---
$ws.Cells.Item(163, 18).Clear()

# --- Append rows 164-176 (new question records) ---
# 'ano' (column D) is stored as text throughout rows 27-163; pre-format
# the new rows' D cells as Text so typing a year-looking value keeps that.
$ws.Range("D164:D176").NumberFormat = "@"

# Row 164
$ws.Cells.Item(164, 1).Value2 = 164
$ws.Cells.Item(164, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(164, 3).Value2 = "BANCO DO BRASIL"
$ws.Cells.Item(164, 4).Value2 = "2021"
$ws.Cells.Item(164, 5).Value2 = "De acordo com a norma-padrão da língua portuguesa, o emprego adequado da vírgula está plenamente atendido em:"
$ws.Cells.Item(164, 6).Value2 = "Português"
$ws.Cells.Item(164, 7).Value2 = "Vírgula"
$ws.Cells.Item(164, 8).Value2 = "Médio"
$ws.Cells.Item(164, 9).Value2 = "ME"
$ws.Cells.Item(164, 10).Value2 = "O ensino remoto, com a pandemia de Covid-19 passou a fazer parte do processo de escolarização em todo o Brasil."
$ws.Cells.Item(164, 11).Value2 = "A melhor fase do ensino on-line tem sido vivida, atualmente embora permaneça a dúvida se é possível ensinar às crianças de forma remota."
$ws.Cells.Item(164, 12).Value2 = "Como o país não tinha experiências significativas no ensino remoto, precisou aderir à prática de forma emergencial."
$ws.Cells.Item(164, 13).Value2 = "A qualidade do ensino remoto era questionada, no passado porém o aprendizado conta com tecnologias que garantem ótimos resultados."
$ws.Cells.Item(164, 14).Value2 = "Um grande número de pesquisadores tem procurado avaliar, quais são as vantagens e desvantagens da educação a distância."
$ws.Cells.Item(164, 15).Value2 = "C"
$ws.Cells.Item(164, 16).Value2 = 0
$ws.Cells.Item(164, 17).Value2 = 0

# Row 165
$ws.Cells.Item(165, 1).Value2 = 165
$ws.Cells.Item(165, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(165, 3).Value2 = "UNIRIO"
$ws.Cells.Item(165, 4).Value2 = "2019"
$ws.Cells.Item(165, 5).Value2 = "Considere o trecho: “Diolino bolou então o sistema de atendimento direto aos veículos”. Caso fosse necessário reescrevê-lo empregando alguma vírgula e mantendo o sentido original, o resultado, de acordo com as normas pontuação, seria:"
$ws.Cells.Item(165, 6).Value2 = "Português"
$ws.Cells.Item(165, 7).Value2 = "Vírgula"
$ws.Cells.Item(165, 8).Value2 = "Médio"
$ws.Cells.Item(165, 9).Value2 = "ME"
$ws.Cells.Item(165, 10).Value2 = "Diolino, bolou então o sistema de atendimento direto, aos veículos."
$ws.Cells.Item(165, 11).Value2 = "Diolino bolou então, o sistema, de atendimento direto aos veículos."
$ws.Cells.Item(165, 12).Value2 = "Diolino bolou então o sistema, de atendimento direto aos veículos."
$ws.Cells.Item(165, 13).Value2 = "Diolino bolou, então, o sistema de atendimento direto aos veículos."
$ws.Cells.Item(165, 14).Value2 = "Diolino bolou, então o sistema de atendimento direto aos veículos."
$ws.Cells.Item(165, 15).Value2 = "D"
$ws.Cells.Item(165, 16).Value2 = 0
$ws.Cells.Item(165, 17).Value2 = 0

# Row 166
$ws.Cells.Item(166, 1).Value2 = 166
$ws.Cells.Item(166, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(166, 3).Value2 = "UNIRIO"
$ws.Cells.Item(166, 4).Value2 = "2019"
$ws.Cells.Item(166, 5).Value2 = "A frase em que a vírgula está empregada adequadamente é:"
$ws.Cells.Item(166, 6).Value2 = "Português"
$ws.Cells.Item(166, 7).Value2 = "Vírgula"
$ws.Cells.Item(166, 8).Value2 = "Médio"
$ws.Cells.Item(166, 9).Value2 = "ME"
$ws.Cells.Item(166, 10).Value2 = "A tela do computador, é a janela que descortina o mundo."
$ws.Cells.Item(166, 11).Value2 = "O investimento deve ser feito na área que, pode salvar vidas."
$ws.Cells.Item(166, 12).Value2 = "A vaga é para programador, que tem salário acima da média."
$ws.Cells.Item(166, 13).Value2 = "Concluíram, que não há mais como parar o avanço tecnológico."
$ws.Cells.Item(166, 14).Value2 = "É muito importante, que os investimentos na área tecnológica continuem."
$ws.Cells.Item(166, 15).Value2 = "C"
$ws.Cells.Item(166, 16).Value2 = 0
$ws.Cells.Item(166, 17).Value2 = 0

# Row 167
$ws.Cells.Item(167, 1).Value2 = 167
$ws.Cells.Item(167, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(167, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(167, 4).Value2 = "2018"
$ws.Cells.Item(167, 5).Value2 = "Considere a seguinte passagem: “Dentro do amplo território portuário, os planejadores urbanos que idealizaram o Plano Porto do Rio haviam concentrado investimentos simbólicos e materiais nos arredores da praça Mauá, situada na convergência do bairro da Saúde com a avenida Rio Branco”.`r`nA reescritura que mantém os aspectos informacionais do trecho e respeita as normas de emprego dos sinais de pontuação é a seguinte:"
$ws.Cells.Item(167, 6).Value2 = "Português"
$ws.Cells.Item(167, 7).Value2 = "Vírgula"
$ws.Cells.Item(167, 8).Value2 = "Médio"
$ws.Cells.Item(167, 9).Value2 = "ME"
$ws.Cells.Item(167, 10).Value2 = "Os planejadores urbanos, que idealizaram dentro do amplo território portuário o Plano Porto do Rio haviam concentrado investimentos simbólicos e materiais nos arredores da praça Mauá, situada na convergência do bairro da Saúde com a avenida Rio Branco."
$ws.Cells.Item(167, 11).Value2 = "Dentro do amplo território portuário, os planejadores urbanos que idealizaram o Plano Porto do Rio, haviam concentrado investimentos simbólicos e materiais nos arredores da praça Mauá, situada na convergência do bairro da Saúde com a avenida Rio Branco."
$ws.Cells.Item(167, 12).Value2 = "Os planejadores urbanos que idealizaram, dentro do amplo território portuário, o Plano Porto do Rio haviam concentrado, investimentos simbólicos e materiais nos arredores da praça Mauá, situada na convergência do bairro da Saúde com a avenida Rio Branco."
$ws.Cells.Item(167, 13).Value2 = "Os planejadores urbanos que idealizaram, dentro do amplo território portuário, o Plano Porto do Rio haviam concentrado investimentos simbólicos e materiais nos arredores da praça Mauá, situada na convergência do bairro da Saúde com a avenida Rio Branco."
$ws.Cells.Item(167, 14).Value2 = "Dentro do amplo, território portuário, os planejadores urbanos que idealizaram o Plano Porto do Rio haviam concentrado investimentos simbólicos e materiais nos arredores da praça Mauá situada na convergência do bairro da Saúde com a avenida Rio Branco."
$ws.Cells.Item(167, 15).Value2 = "D"
$ws.Cells.Item(167, 16).Value2 = 0
$ws.Cells.Item(167, 17).Value2 = 0

# Row 168
$ws.Cells.Item(168, 1).Value2 = 168
$ws.Cells.Item(168, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(168, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(168, 4).Value2 = "2018"
$ws.Cells.Item(168, 5).Value2 = "Considere-se o emprego da primeira vírgula no trecho transcrito abaixo. “Naquele momento, quem o visse de perto perceberia o suor escorrendo frio por seu rosto”. A vírgula é empregada pelo mesmo motivo em:"
$ws.Cells.Item(168, 6).Value2 = "Português"
$ws.Cells.Item(168, 7).Value2 = "Vírgula"
$ws.Cells.Item(168, 8).Value2 = "Médio"
$ws.Cells.Item(168, 9).Value2 = "ME"
$ws.Cells.Item(168, 10).Value2 = "A falta não foi dentro da área, mas o juiz deu pênalti."
$ws.Cells.Item(168, 11).Value2 = "O atacante passou pelo zagueiro, passou pelo goleiro e fez o gol."
$ws.Cells.Item(168, 12).Value2 = "Lúcio, atrase a bola para o goleiro!"
$ws.Cells.Item(168, 13).Value2 = "O juiz verificou as balizas, a bola, as marcações do campo e deu início à partida."
$ws.Cells.Item(168, 14).Value2 = "No campo de jogo, Lúcio se sentia livre."
$ws.Cells.Item(168, 15).Value2 = "E"
$ws.Cells.Item(168, 16).Value2 = 0
$ws.Cells.Item(168, 17).Value2 = 0

# Row 169
$ws.Cells.Item(169, 1).Value2 = 169
$ws.Cells.Item(169, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(169, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(169, 4).Value2 = "2018"
$ws.Cells.Item(169, 5).Value2 = "A vírgula está empregada corretamente em:"
$ws.Cells.Item(169, 6).Value2 = "Português"
$ws.Cells.Item(169, 7).Value2 = "Vírgula"
$ws.Cells.Item(169, 8).Value2 = "Médio"
$ws.Cells.Item(169, 9).Value2 = "ME"
$ws.Cells.Item(169, 10).Value2 = "A divulgação de histórias falsas pode ter consequências reais desastrosas: prejuízos, financeiros e constrangimentos às empresas."
$ws.Cells.Item(169, 11).Value2 = "As novas tecnologias, criaram um abismo ao separar quem está conectado de quem não faz parte do mundo digital."
$ws.Cells.Item(169, 12).Value2 = "As pessoas tendem a aceitar apenas as declarações que confirmam aquilo que corresponde, às suas crenças."
$ws.Cells.Item(169, 13).Value2 = "Os jornalistas devem verificar as fontes citadas, cruzar dados e checar se as informações refletem a realidade."
$ws.Cells.Item(169, 14).Value2 = "Os consumidores de notícias não agem como cientistas porque não estão preocupados em conferir, pontos de vista alternativos."
$ws.Cells.Item(169, 15).Value2 = "D"
$ws.Cells.Item(169, 16).Value2 = 0
$ws.Cells.Item(169, 17).Value2 = 0

# Row 170
$ws.Cells.Item(170, 1).Value2 = 170
$ws.Cells.Item(170, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(170, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(170, 4).Value2 = "2018"
$ws.Cells.Item(170, 5).Value2 = "A vírgula está empregada de acordo com a norma-padrão da língua portuguesa em:"
$ws.Cells.Item(170, 6).Value2 = "Português"
$ws.Cells.Item(170, 7).Value2 = "Vírgula"
$ws.Cells.Item(170, 8).Value2 = "Médio"
$ws.Cells.Item(170, 9).Value2 = "ME"
$ws.Cells.Item(170, 10).Value2 = "A acessibilidade é a possibilidade que as pessoas, têm de atingir o destino desejado."
$ws.Cells.Item(170, 11).Value2 = "A mobilidade urbana tem, forte impacto, sobre o espaço e os recursos naturais."
$ws.Cells.Item(170, 12).Value2 = "As políticas públicas, devem priorizar os meios de transporte coletivo, nas cidades."
$ws.Cells.Item(170, 13).Value2 = "Como alertam os pesquisadores, é preciso discutir estratégias de mobilidade urbana."
$ws.Cells.Item(170, 14).Value2 = "Nos últimos anos aumentou, a insatisfação das pessoas com os engarrafamentos."
$ws.Cells.Item(170, 15).Value2 = "D"
$ws.Cells.Item(170, 16).Value2 = 0
$ws.Cells.Item(170, 17).Value2 = 0

# Row 171
$ws.Cells.Item(171, 1).Value2 = 171
$ws.Cells.Item(171, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(171, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(171, 4).Value2 = "2018"
$ws.Cells.Item(171, 5).Value2 = "Considere o trecho “Num mundo em que as notícias são plantadas pela internet, em que muitos sites servem a qualquer mentira.”.`r`nA única reescritura que, além de conservar o conteúdo informacional, emprega os sinais de pontuação de acordo com a norma-padrão é:"
$ws.Cells.Item(171, 6).Value2 = "Português"
$ws.Cells.Item(171, 7).Value2 = "Vírgula"
$ws.Cells.Item(171, 8).Value2 = "Médio"
$ws.Cells.Item(171, 9).Value2 = "ME"
$ws.Cells.Item(171, 10).Value2 = "Num mundo em que as notícias, são plantadas pela internet, em que muitos sites servem a qualquer mentira."
$ws.Cells.Item(171, 11).Value2 = "Num mundo em que muitos sites servem a qualquer mentira, em que as notícias são plantadas pela internet."
$ws.Cells.Item(171, 12).Value2 = "Num mundo em que, pela internet, as notícias são plantadas em que muitos sites, servem a qualquer mentira."
$ws.Cells.Item(171, 13).Value2 = "Num mundo, em que as notícias são plantadas pela internet em muitos sites que servem a qualquer mentira."
$ws.Cells.Item(171, 14).Value2 = "Num mundo em que, as notícias são plantadas pela internet e em que, muitos sites servem a qualquer mentira."
$ws.Cells.Item(171, 15).Value2 = "B"
$ws.Cells.Item(171, 16).Value2 = 0
$ws.Cells.Item(171, 17).Value2 = 0

# Row 172
$ws.Cells.Item(172, 1).Value2 = 172
$ws.Cells.Item(172, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(172, 3).Value2 = "PETROBRAS"
$ws.Cells.Item(172, 4).Value2 = "2018"
$ws.Cells.Item(172, 5).Value2 = "Observe atentamente o uso dos sinais de pontuação do trecho abaixo: “Há, de igual forma, entre os mais afortunados, aqueles ligados à indústria, voltados para a construção civil, o mobiliário, a ourivesaria e o fabrico de bebidas.” Qual das reescrituras desse trecho emprega corretamente os sinais de pontuação?"
$ws.Cells.Item(172, 6).Value2 = "Português"
$ws.Cells.Item(172, 7).Value2 = "Vírgula"
$ws.Cells.Item(172, 8).Value2 = "Médio"
$ws.Cells.Item(172, 9).Value2 = "ME"
$ws.Cells.Item(172, 10).Value2 = "Há, entre os mais afortunados de igual forma, aqueles ligados à indústria voltados para a construção civil, o mobiliário, a ourivesaria, e o fabrico de bebidas."
$ws.Cells.Item(172, 11).Value2 = "De igual forma, há, entre os mais afortunados, aqueles ligados à indústria, voltados para a construção civil, o mobiliário, a ourivesaria e o fabrico de bebidas."
$ws.Cells.Item(172, 12).Value2 = "Entre os mais afortunados, há de igual forma, aqueles ligados à indústria, voltados para a construção civil, o mobiliário, a ourivesaria, e o fabrico de bebidas."
$ws.Cells.Item(172, 13).Value2 = "Há entre os mais afortunados de igual forma, aqueles ligados à indústria, voltados para a construção civil, o mobiliário, a ourivesaria e o fabrico de bebidas."
$ws.Cells.Item(172, 14).Value2 = "De igual forma, entre os mais afortunados, há, aqueles, ligados à indústria, voltados para a construção civil, o mobiliário, a ourivesaria e o fabrico de bebidas."
$ws.Cells.Item(172, 15).Value2 = "B"
$ws.Cells.Item(172, 16).Value2 = 0
$ws.Cells.Item(172, 17).Value2 = 0

# Row 173
$ws.Cells.Item(173, 1).Value2 = 173
$ws.Cells.Item(173, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(173, 3).Value2 = "PETROBRAS"
$ws.Cells.Item(173, 4).Value2 = "2018"
$ws.Cells.Item(173, 5).Value2 = "A vírgula está empregada corretamente em:"
$ws.Cells.Item(173, 6).Value2 = "Português"
$ws.Cells.Item(173, 7).Value2 = "Vírgula"
$ws.Cells.Item(173, 8).Value2 = "Médio"
$ws.Cells.Item(173, 9).Value2 = "ME"
$ws.Cells.Item(173, 10).Value2 = "As grandes metrópoles que se destacaram no apoio à sustentabilidade, foram premiadas pelo mundo inteiro."
$ws.Cells.Item(173, 11).Value2 = "É preciso que futuramente, as cidades tenham melhores condições de vida: habitação, alimentação, saúde, emprego, transporte, educação."
$ws.Cells.Item(173, 12).Value2 = "Não é só o território que acelera o seu processo de urbanização, mas é a própria sociedade brasileira que se transforma cada vez mais em urbana."
$ws.Cells.Item(173, 13).Value2 = "Os estados que possuem os menores percentuais de população vivendo em áreas urbanas, estão concentrados nas regiões Norte e Nordeste."
$ws.Cells.Item(173, 14).Value2 = "Os passageiros, que dependem do transporte coletivo esperam que o futuro lhes ofereça mais comodidade do que o presente."
$ws.Cells.Item(173, 15).Value2 = "C"
$ws.Cells.Item(173, 16).Value2 = 0
$ws.Cells.Item(173, 17).Value2 = 0

# Row 174
$ws.Cells.Item(174, 1).Value2 = 174
$ws.Cells.Item(174, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(174, 3).Value2 = "PETROBRAS"
$ws.Cells.Item(174, 4).Value2 = "2018"
$ws.Cells.Item(174, 5).Value2 = "A vírgula foi plenamente empregada de acordo com as exigências da norma-padrão da língua portuguesa em:"
$ws.Cells.Item(174, 6).Value2 = "Português"
$ws.Cells.Item(174, 7).Value2 = "Vírgula"
$ws.Cells.Item(174, 8).Value2 = "Médio"
$ws.Cells.Item(174, 9).Value2 = "ME"
$ws.Cells.Item(174, 10).Value2 = "A conexão é feita por meio de uma plataforma específica, e os conteúdos, podem ser acessados pelos dispositivos móveis dos passageiros."
$ws.Cells.Item(174, 11).Value2 = "O mercado brasileiro de automóveis, ainda é muito grande, porém não é capaz de absorver uma presença maior de produtos vindos do exterior."
$ws.Cells.Item(174, 12).Value2 = "Depois de chegarem às telas dos computadores e celulares, as notícias estarão disponíveis em voos internacionais."
$ws.Cells.Item(174, 13).Value2 = "Os últimos dados mostram que, muitas economias apresentam crescimento e inflação baixa, fazendo com que os juros cresçam pouco."
$ws.Cells.Item(174, 14).Value2 = "Pode ser que haja uma grande procura de carros importados, mas as montadoras vão fazer os cálculos e ver, se a importação vale a pena."
$ws.Cells.Item(174, 15).Value2 = "C"
$ws.Cells.Item(174, 16).Value2 = 0
$ws.Cells.Item(174, 17).Value2 = 0

# Row 175
$ws.Cells.Item(175, 1).Value2 = 175
$ws.Cells.Item(175, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(175, 3).Value2 = "BANCO DO BRASIL"
$ws.Cells.Item(175, 4).Value2 = "2021"
$ws.Cells.Item(175, 5).Value2 = "O período em que o sinal de dois pontos é empregado para introduzir uma enumeração, como no trecho que segue “demanda garantida” (parágrafo 2), é:"
$ws.Cells.Item(175, 6).Value2 = "Português"
$ws.Cells.Item(175, 7).Value2 = "Dois"
$ws.Cells.Item(175, 8).Value2 = "Médio"
$ws.Cells.Item(175, 9).Value2 = "ME"
$ws.Cells.Item(175, 10).Value2 = "A remuneração faz parte do conjunto de ganhos de um prestador de serviço; ou seja: todos os ganhos auferidos pela pessoa compõem sua remuneração."
$ws.Cells.Item(175, 11).Value2 = "As horas extras, o vale-transporte e o plano de saúde podem fazer parte da remuneração: muitos trabalhadores escolhem seus empregos com base nessas vantagens."
$ws.Cells.Item(175, 12).Value2 = "O gerente informou aos candidatos como seria a remuneração pelos serviços: “O valor mensal vai depender de diversos itens, a serem combinados.”"
$ws.Cells.Item(175, 13).Value2 = "Muitos itens já fizeram papel de dinheiro: o sal, usado até hoje por tribos da Etiópia, a cachaça, utilizada no Brasil colonial, e o bacalhau, antes usado na Escandinávia."
$ws.Cells.Item(175, 14).Value2 = "O tabaco também já foi usado como moeda de troca: no século XVIII, o estado americano de Virginia adotou esse método."
$ws.Cells.Item(175, 15).Value2 = "D"
$ws.Cells.Item(175, 16).Value2 = 0
$ws.Cells.Item(175, 17).Value2 = 0
$ws.Cells.Item(175, 18).Value2 = "6833f320-6441-4318-89ab-c90bb9f149e6.png"

# Row 176
$ws.Cells.Item(176, 1).Value2 = 176
$ws.Cells.Item(176, 2).Value2 = "CESGRANRIO"
$ws.Cells.Item(176, 3).Value2 = "LIQUIGÁS"
$ws.Cells.Item(176, 4).Value2 = "2015"
$ws.Cells.Item(176, 5).Value2 = "Em “Posso intensamente desejar que o problema mais urgente se resolva: o da fome.`", os dois-pontos cumprem o papel de introduzir uma"
$ws.Cells.Item(176, 6).Value2 = "Português"
$ws.Cells.Item(176, 7).Value2 = "Dois"
$ws.Cells.Item(176, 8).Value2 = "Médio"
$ws.Cells.Item(176, 9).Value2 = "ME"
$ws.Cells.Item(176, 10).Value2 = "enumeração"
$ws.Cells.Item(176, 11).Value2 = "explanação"
$ws.Cells.Item(176, 12).Value2 = "retificação"
$ws.Cells.Item(176, 13).Value2 = "especificação"
$ws.Cells.Item(176, 14).Value2 = "contradição."
$ws.Cells.Item(176, 15).Value2 = "D"
$ws.Cells.Item(176, 16).Value2 = 0
$ws.Cells.Item(176, 17).Value2 = 0

# --- Refresh the declared dimension to cover the appended rows ---
$ws.UsedRange | Out-Null
